$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data scraped (betexplorer-style) to append to the sheet.
$rows = @(
    @{
        Row = 32
        Idx = 31
        Pais = "kuwait"
        Torneio = "premier-league"
        Temporada = "2023-2024"
        Data = 45232.64583333334
        Home = "Al Naser"
        HomeGols = 3
        Away = "Al Arabi"
        AwayGols = 1
        HomeOpenOdds = 3.81
        HomeOpenData = "01/11/2023 03:42"
        HomeCloseOdds = 4.23
        HomeCloseData = "02/11/2023 15:05"
        DrawOpenOdds = 3.78
        DrawOpenData = "01/11/2023 03:42"
        DrawCloseOdds = 3.85
        DrawCloseData = "02/11/2023 15:08"
        AwayOpenOdds = 1.66
        AwayOpenData = "01/11/2023 03:42"
        AwayCloseOdds = 1.7
        AwayCloseData = "02/11/2023 15:05"
        Url = "https://www.betexplorer.com/football/kuwait/premier-league/al-naser-al-arabi-kuwait/IBBEGPPH/"
    },
    @{
        Row = 33
        Idx = 32
        Pais = "kuwait"
        Torneio = "premier-league"
        Temporada = "2023-2024"
        Data = 45232.76388888889
        Home = "Al Kuwait"
        HomeGols = 4
        Away = "Al-Fahaheel"
        AwayGols = 1
        HomeOpenOdds = 1.24
        HomeOpenData = "01/11/2023 08:42"
        HomeCloseOdds = 1.18
        HomeCloseData = "02/11/2023 18:17"
        DrawOpenOdds = 5.38
        DrawOpenData = "01/11/2023 08:42"
        DrawCloseOdds = 6.72
        DrawCloseData = "02/11/2023 18:17"
        AwayOpenOdds = 7.32
        AwayOpenData = "01/11/2023 08:42"
        AwayCloseOdds = 11.42
        AwayCloseData = "02/11/2023 18:17"
        Url = "https://www.betexplorer.com/football/kuwait/premier-league/al-kuwait-al-fahaheel/QyBIFquO/"
    }
)

# The last pre-existing data row (30 => sheet row 31) carries the exact
# cell styles (bold/bordered index cell, formatted date cell) that new
# rows must reuse. Copy its formatting down before writing the new values.
$templateRow = 31

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Range("A$templateRow").Copy() | Out-Null
    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null

    $ws.Range("E$templateRow").Copy() | Out-Null
    $ws.Range("E$row").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 1).Value = $r.Idx

    $ws.Cells.Item($row, 2).Value = $r.Pais
    $ws.Cells.Item($row, 3).Value = $r.Torneio
    $ws.Cells.Item($row, 4).Value = $r.Temporada

    $ws.Cells.Item($row, 5).Value = $r.Data

    $ws.Cells.Item($row, 6).Value = $r.Home
    $ws.Cells.Item($row, 7).Value = $r.HomeGols
    $ws.Cells.Item($row, 8).Value = $r.Away
    $ws.Cells.Item($row, 9).Value = $r.AwayGols

    $ws.Cells.Item($row, 10).Value = $r.HomeOpenOdds
    $ws.Cells.Item($row, 11).Value = $r.HomeOpenData
    $ws.Cells.Item($row, 12).Value = $r.HomeCloseOdds
    $ws.Cells.Item($row, 13).Value = $r.HomeCloseData

    $ws.Cells.Item($row, 14).Value = $r.DrawOpenOdds
    $ws.Cells.Item($row, 15).Value = $r.DrawOpenData
    $ws.Cells.Item($row, 16).Value = $r.DrawCloseOdds
    $ws.Cells.Item($row, 17).Value = $r.DrawCloseData

    $ws.Cells.Item($row, 18).Value = $r.AwayOpenOdds
    $ws.Cells.Item($row, 19).Value = $r.AwayOpenData
    $ws.Cells.Item($row, 20).Value = $r.AwayCloseOdds
    $ws.Cells.Item($row, 21).Value = $r.AwayCloseData

    $ws.Cells.Item($row, 22).Value = $r.Url
}

$excel.CutCopyMode = $false
